$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.375.26'
$ws.Range('E2').Value = '  -7.30%  '
$ws.Range('D3').Value = '2.867.30'
$ws.Range('E3').Value = '  -10.51%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '472.82'
$ws.Range('E5').Value = '  -11.53%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '125.44'
$ws.Range('E6').Value = '  -6.63%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '2.863.62'
$ws.Range('E8').Value = '  -10.72%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.402'
$ws.Range('E9').Value = '  -11.81%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.64'
$ws.Range('E10').Value = '  -11.79%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0960'
$ws.Range('E11').Value = '  -15.51%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.329'
$ws.Range('E12').Value = '  -15.66%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.122'
$ws.Range('E13').Value = '  -4.35%  '
$ws.Range('D14').Value = '3.359.91'
$ws.Range('E14').Value = '  -10.69%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '22.86'
$ws.Range('E15').Value = '  -11.31%  '
$ws.Range('D16').Value = '54.347.31'
$ws.Range('E16').Value = '  -7.44%  '
$ws.Range('D17').Value = '2.874.29'
$ws.Range('E17').Value = '  -10.51%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0000134'
$ws.Range('E18').Value = '  -14.73%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.27'
$ws.Range('E19').Value = '  -10.65%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.48'
$ws.Range('E20').Value = '  -13.00%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.02'
$ws.Range('E21').Value = '  -14.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '303.18'
$ws.Range('E22').Value = '  -15.72%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.443'
$ws.Range('E24').Value = '  -14.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '58.79'
$ws.Range('E25').Value = '  -16.02%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('E27').Value = '  -9.90%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('D29').Value = '0.0₃0804'
$ws.Range('E29').Value = '  -15.70%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.24'
$ws.Range('E30').Value = '  -11.58%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.12'
$ws.Range('E31').Value = '  -6.58%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.14'
$ws.Range('E32').Value = '  -12.70%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.96'
$ws.Range('E33').Value = '  -12.53%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.61'
$ws.Range('E34').Value = '  -16.23%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.19'
$ws.Range('E35').Value = '  -14.20%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '134.07'
$ws.Range('E36').Value = '  -16.85%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.40'
$ws.Range('E37').Value = '  -14.97%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.20'
$ws.Range('E38').Value = '  -16.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '22.58'
$ws.Range('E39').Value = '  -12.27%  '
$ws.Range('D40').Value = '2.895.02'
$ws.Range('E40').Value = '  -10.53%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0609'
$ws.Range('E41').Value = '  -13.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '35.63'
$ws.Range('E43').Value = '  -13.17%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.964'
$ws.Range('E44').Value = '  -11.72%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.598'
$ws.Range('E45').Value = '  -16.16%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.38'
$ws.Range('E46').Value = '  -15.37%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.30'
$ws.Range('E47').Value = '  -12.52%  '
$ws.Range('D48').Value = '2.038.01'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.28'
$ws.Range('E49').Value = '  -15.66%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.68'
$ws.Range('E50').Value = '  -14.70%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0210'
$ws.Range('E51').Value = '  -12.09%  '
